# The deck has three tables (on slides 14, 15 and 16) that were still using
# the deck's custom "Table_0" style. Re-apply a different built-in table
# style (the one PowerPoint's Table Styles gallery identifies by the GUID
# below) to each of them, matching the table style id written into the
# slide XML for all three tables.

$p = $ppt.ActivePresentation

$newStyleId = "{F6E1076A-A33A-4C04-B0B4-65FB0608403F}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    $tableShape = $slide.Shapes.Item(1)
    if ($tableShape.HasTable) {
        $tableShape.Table.ApplyStyle($newStyleId)
    }
}
